# Cardiac arrhythmias - additional genes: split panel metadata into its own tab.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Refresh the panel_query_time timestamps on the "data" sheet ----------
$ws.Range("F2").Value = "2021-10-05 14:19:24.597084"
$ws.Range("F3").Value = "2021-10-05 14:19:24.597092"

# --- Add the new "metadata" worksheet right after "data" ------------------
# Duplicate "data" (so sheetPr/outline/page-setup carry over identically)
# then wipe its contents before repopulating with the metadata columns.
$ws.Copy($null, $ws)
$meta = $wb.Worksheets.Item("data (2)")
$meta.Name = "metadata"
$meta.Cells.Clear()

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("B1:G1").Style = $ws.Range("B1").Style

# Data row
$meta.Range("A2").Value = 0
$meta.Range("A2").Style = $ws.Range("A2").Style

$meta.Range("B2").Value = "Cardiac arrhythmias - additional genes"

$meta.Range("C2").Value = 843

# data_version must stay textual ("1.12"), not become the number 1.12
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.12"
$meta.Range("D2").NumberFormat = "General"
$meta.Range("D2").Style = $ws.Range("B2").Style

$meta.Range("E2").Value = "2021-08-10T08:14:00.626794Z"
$meta.Range("F2").Value = "2021-10-05 14:19:24.593352"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/843/?format=json"

# Keep "data" as the active/selected sheet, as it was before this edit.
$ws.Activate()
